$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for inorganic carbon data columns (UD, LD) next to existing ones
$ws.Range("G1").Value = "UD"
$ws.Range("H1").Value = "LD"

# Update selection to reflect the new active cell after the edit (next empty column)
$ws.Range("I1").Select()
